$wb = $excel.ActiveWorkbook

# --- Sheet 2: matematyka stosowana_1_1 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("C2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "Analiza_matematyczna_I_lecture_1"
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("B4").Value = "Analiza_matematyczna_I_practicals_2_grp_1"
$ws.Range("C4").Value = "Algebra_liniowa_z_geometrią_analityczną_I_practicals_1_grp_1"
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = "Wstęp_do_logiki_i_teorii_mnogości_practicals_1_grp_1"
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = "Algebra_liniowa_z_geometrią_analityczną_I_lecture_1"
$ws.Range("F5").Value = "Wstęp_do_logiki_i_teorii_mnogości_lecture_1"
$ws.Range("D6").ClearContents()
$ws.Range("B7").Value = "Analiza_matematyczna_I_practicals_1_grp_1"
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = "Wstęp_do_obliczeń_symbolicznych_laboratories_1_grp_1"
$ws.Range("E7").Value = "Analiza_matematyczna_I_lecture_2"
$ws.Range("F7").Value = "Technologie_informatyczne_I_laboratories_1_grp_1"

# --- Sheet 3: matematyka stosowana_1_2 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = "Wstęp_do_obliczeń_symbolicznych_laboratories_1_grp_2"
$ws.Range("C2").ClearContents()
$ws.Range("C3").Value = "Analiza_matematyczna_I_lecture_1"
$ws.Range("D3").Value = "Analiza_matematyczna_I_practicals_2_grp_2"
$ws.Range("E4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").Value = "Algebra_liniowa_z_geometrią_analityczną_I_lecture_1"
$ws.Range("F5").Value = "Wstęp_do_logiki_i_teorii_mnogości_lecture_1"
$ws.Range("C6").Value = "Analiza_matematyczna_I_practicals_1_grp_2"
$ws.Range("D6").Value = "Algebra_liniowa_z_geometrią_analityczną_I_practicals_1_grp_2"
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = "Wstęp_do_logiki_i_teorii_mnogości_practicals_1_grp_2"
$ws.Range("C7").Value = "Technologie_informatyczne_I_laboratories_1_grp_2"
$ws.Range("E7").Value = "Analiza_matematyczna_I_lecture_2"

# --- Sheet 4: matematyka stosowana_2_1 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = "Analiza_matematyczna_II_practicals_2_grp_1"
$ws.Range("F2").Value = "Komputerowe_obliczenia_matematyczne_laboratories_1_grp_1"
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "Analiza_matematyczna_II_lecture_1"
$ws.Range("E3").Value = "Analiza_matematyczna_II_practicals_1_grp_1"
$ws.Range("F3").Value = "Podstawy_probabilistyki_lecture_1"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = "Podstawy_ekonomii_i_przedsiębiorczości_lecture_1"
$ws.Range("F4").ClearContents()
$ws.Range("B5").Value = "Technologie_informatyczne_II_laboratories_1_grp_1"
$ws.Range("C5").Value = "Algebra_liniowa_z_geometrią_analityczną_II_practicals_1_grp_1"
$ws.Range("D5").Value = "Topologia_przestrzeni_metrycznych_practicals_1_grp_1"
$ws.Range("E5").Value = "Algebra_liniowa_z_geometrią_analityczną_II_lecture_1"
$ws.Range("F6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "Topologia_przestrzeni_metrycznych_lecture_1"
$ws.Range("F7").Value = "Podstawy_probabilistyki_practicals_1_grp_1"

# --- Sheet 5: matematyka stosowana_2_2 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("C2").Value = "Podstawy_probabilistyki_practicals_1_grp_2"
$ws.Range("F2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "Analiza_matematyczna_II_lecture_1"
$ws.Range("D3").Value = "Topologia_przestrzeni_metrycznych_practicals_1_grp_2"
$ws.Range("E3").Value = "Technologie_informatyczne_II_laboratories_1_grp_2"
$ws.Range("F3").Value = "Podstawy_probabilistyki_lecture_1"
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").Value = "Podstawy_ekonomii_i_przedsiębiorczości_lecture_1"
$ws.Range("F4").Value = "Komputerowe_obliczenia_matematyczne_laboratories_1_grp_2"
$ws.Range("C5").Value = "Analiza_matematyczna_II_practicals_2_grp_2"
$ws.Range("E5").Value = "Algebra_liniowa_z_geometrią_analityczną_II_lecture_1"
$ws.Range("F5").Value = "Algebra_liniowa_z_geometrią_analityczną_II_practicals_1_grp_2"
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "Topologia_przestrzeni_metrycznych_lecture_1"
$ws.Range("D7").ClearContents()
$ws.Range("F7").Value = "Analiza_matematyczna_II_practicals_1_grp_2"

